$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '304.43'
Set-TextValue 'E2' '0.00%'
Set-TextValue 'D3' '35.66'
Set-TextValue 'E3' '-3.81%'
Set-TextValue 'D4' '5.054'
Set-TextValue 'E4' '0.94%'
Set-TextValue 'D5' '0.07880'
Set-TextValue 'E5' '0.05%'
Set-TextValue 'D6' '2.114'
Set-TextValue 'E6' '-4.20%'
Set-TextValue 'B7' 'GateToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D7' '4.125'
Set-TextValue 'E7' '2.61%'
Set-TextValue 'B8' 'KuCoinToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue 'D8' '7.924'
Set-TextValue 'E8' '-1.12%'
Set-TextValue 'B9' 'MXToken'
Set-TextValue 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.9242'
Set-TextValue 'E9' '0.40%'
Set-TextValue 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.09775'
Set-TextValue 'E10' '1.39%'
Set-TextValue 'B11' 'WazirX'
Set-TextValue 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1835'
Set-TextValue 'E11' '-2.90%'
Set-TextValue 'B12' 'MandalaExchangeToken'
Set-TextValue 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.08652'
Set-TextValue 'E12' '0.81%'
Set-TextValue 'B13' 'BitrueCoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03535'
Set-TextValue 'E13' '-4.09%'
Set-TextValue 'B14' 'BitMartToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09904'
Set-TextValue 'E14' '-0.75%'
Set-TextValue 'B15' 'BitForexToken'
Set-TextValue 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001432'
Set-TextValue 'E15' '-4.16%'
Set-TextValue 'B16' 'TigerCash'
Set-TextValue 'C16' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D16' '0.005694'
Set-TextValue 'E16' '0.09%'
Set-TextValue 'B17' 'LEO'
Set-TextValue 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D17' '3.454'
Set-TextValue 'E17' '-0.43%'
Set-TextValue 'D18' '2.640'
Set-TextValue 'E18' '17.37%'
Set-TextValue 'D19' '0.3371'
Set-TextValue 'E19' '-1.27%'
Set-TextValue 'D20' '0.1339'
Set-TextValue 'E20' '1.66%'
Set-TextValue 'D21' '5.179'
Set-TextValue 'E21' '9.01%'
Set-TextValue 'D22' '0.2212'
Set-TextValue 'E22' '0.53%'
Set-TextValue 'D23' '0.04497'
Set-TextValue 'E23' '-1.23%'
Set-TextValue 'D24' '0.001236'
Set-TextValue 'E24' '0.15%'
Set-TextValue 'E25' '8.70%'
Set-TextValue 'D26' '0.0001304'
Set-TextValue 'E26' '-6.78%'
Set-TextValue 'D27' '0.0004758'
Set-TextValue 'E27' '0.14%'
Set-TextValue 'D39' '0.01827'
Set-TextValue 'E39' '-0.78%'
Set-TextValue 'D40' '0.04702'
Set-TextValue 'D41' '0.007867'
Set-TextValue 'E41' '-3.17%'
Set-TextValue 'D42' '0.1386'
Set-TextValue 'E42' '-0.87%'
Set-TextValue 'D43' '0.007746'
Set-TextValue 'E43' '2.58%'
Set-TextValue 'D44' '0.002195'
Set-TextValue 'E44' '-1.47%'
Set-TextValue 'D45' '0.01116'
Set-TextValue 'E45' '6.31%'
Set-TextValue 'D46' '0.00006279'
Set-TextValue 'E46' '0.15%'
Set-TextValue 'D47' '0.00000000751'
Set-TextValue 'E47' '0.28%'
Set-TextValue 'E48' '0.21%'
Set-TextValue 'D49' '50.63'
Set-TextValue 'E49' '69.63%'
Set-TextValue 'E50' '10.61%'
Set-TextValue 'D51' '0.00002104'
Set-TextValue 'E51' '0.28%'
